$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 621.2857
$ws.Range("I5").Value = 638.1539
$ws.Range("K5").Value = 638.1539
$ws.Range("M5").Value = -523.1539
$ws.Range("H8").Value = 64.666664
$ws.Range("I8").Value = 69.59999999999999
$ws.Range("J8").Value = 40
$ws.Range("K8").Value = 208.8
$ws.Range("L8").Value = 120
$ws.Range("M8").Value = -69.79999999999998
$ws.Range("N8").Value = -398
$ws.Range("H9").Value = 253.85715
$ws.Range("I9").Value = 257.86365
$ws.Range("J9").Value = 239.16667
$ws.Range("K9").Value = 257.86365
$ws.Range("L9").Value = 239.16667
$ws.Range("M9").Value = -88.86365000000001
$ws.Range("N9").Value = -577.1666700000001
$ws.Range("H29").Value = 837.125
$ws.Range("I29").Value = 232.83333
$ws.Range("J29").Value = 2650
$ws.Range("K29").Value = 698.49999
$ws.Range("L29").Value = 7950
$ws.Range("M29").Value = -417.49999
$ws.Range("N29").Value = -8512
$ws.Range("H34").Value = 14656.5
$ws.Range("I34").Value = 14656.5
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 14656.5
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -14453.5
$ws.Range("N34").Value = $null
$ws.Range("H36").Value = 14656.5
$ws.Range("I36").Value = 14656.5
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 14656.5
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -13941.5
$ws.Range("N36").Value = $null
$ws.Range("H46").Value = 1520
$ws.Range("I46").Value = 1100.5
$ws.Range("J46").Value = 1799.6666
$ws.Range("K46").Value = 3301.5
$ws.Range("L46").Value = 5398.9998
$ws.Range("M46").Value = -3182.5
$ws.Range("N46").Value = -5636.9998
$ws.Range("H60").Value = 1520
$ws.Range("I60").Value = 1100.5
$ws.Range("J60").Value = 1799.6666
$ws.Range("K60").Value = 3301.5
$ws.Range("L60").Value = 5398.9998
$ws.Range("M60").Value = -2817.5
$ws.Range("N60").Value = -6366.9998
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = $null
$ws.Range("H103").Value = 599.5
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 599.5
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 1798.5
$ws.Range("M103").Value = $null
$ws.Range("N103").Value = -2970.5
$ws.Range("H116").Value = 39967.824
$ws.Range("I116").Value = 80694.5
$ws.Range("J116").Value = 3766.3333
$ws.Range("K116").Value = 80694.5
$ws.Range("L116").Value = 3766.3333
$ws.Range("M116").Value = -77252.5
$ws.Range("N116").Value = -10650.3333
$ws.Range("H125").Value = 3828.3333
$ws.Range("I125").Value = 3625
$ws.Range("J125").Value = 3991
$ws.Range("K125").Value = 32625
$ws.Range("L125").Value = 35919
$ws.Range("M125").Value = -30165
$ws.Range("N125").Value = -40839
$ws.Range("H135").Value = 27779038
$ws.Range("I135").Value = 29412982
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 264716838
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -264714303
$ws.Range("N135").Value = -23070
$ws.Range("H137").Value = 2641.875
$ws.Range("I137").Value = 1479.3125
$ws.Range("K137").Value = 4437.9375
$ws.Range("M137").Value = -1887.9375
$ws.Range("H138").Value = 3972.2124
$ws.Range("J138").Value = 4226.4204
$ws.Range("L138").Value = 12679.2612
$ws.Range("N138").Value = -22959.2612
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null
$ws.Range("H141").Value = 1108.0834
$ws.Range("I141").Value = 1108.0834
$ws.Range("K141").Value = 3324.2502
$ws.Range("M141").Value = 1855.7498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2159.4722
$ws.Range("I2").Value = 1742.6
$ws.Range("J2").Value = 2680.5625
$ws.Range("K2").Value = 1742.6
$ws.Range("L2").Value = 2680.5625
$ws.Range("M2").Value = -1629.6
$ws.Range("N2").Value = -2906.5625
$ws.Range("H4").Value = 1042.4286
$ws.Range("I4").Value = 460
$ws.Range("J4").Value = 2498.5
$ws.Range("K4").Value = 460
$ws.Range("L4").Value = 2498.5
$ws.Range("M4").Value = -344
$ws.Range("N4").Value = -2730.5
$ws.Range("H5").Value = 391.7
$ws.Range("I5").Value = 367.14285
$ws.Range("J5").Value = 449
$ws.Range("K5").Value = 367.14285
$ws.Range("L5").Value = 449
$ws.Range("M5").Value = -255.14285
$ws.Range("N5").Value = -673
$ws.Range("H10").Value = 500
$ws.Range("I10").Value = 500
$ws.Range("K10").Value = 500
$ws.Range("M10").Value = -330
$ws.Range("H32").Value = 8203668.5
$ws.Range("I32").Value = 8933103
$ws.Range("J32").Value = 34008
$ws.Range("K32").Value = 8933103
$ws.Range("L32").Value = 34008
$ws.Range("M32").Value = -8932816
$ws.Range("N32").Value = -34582
$ws.Range("H44").Value = 7603750
$ws.Range("J44").Value = 7603750
$ws.Range("L44").Value = 7603750
$ws.Range("N44").Value = -7604726
$ws.Range("H51").Value = 19523.5
$ws.Range("H53").Value = 17021.5
$ws.Range("J53").Value = 17021.5
$ws.Range("L53").Value = 17021.5
$ws.Range("N53").Value = -18385.5
$ws.Range("H63").Value = 5523.5
$ws.Range("I63").Value = 1698
$ws.Range("K63").Value = 1698
$ws.Range("M63").Value = -1012
$ws.Range("H66").Value = 5523.5
$ws.Range("I66").Value = 1698
$ws.Range("K66").Value = 8490
$ws.Range("M66").Value = -5058
$ws.Range("H74").Value = 4632852
$ws.Range("I74").Value = 6946696.5
$ws.Range("K74").Value = 6946696.5
$ws.Range("M74").Value = -6945822.5
$ws.Range("H77").Value = 4632852
$ws.Range("I77").Value = 6946696.5
$ws.Range("K77").Value = 34733482.5
$ws.Range("M77").Value = -34729114.5
$ws.Range("H80").Value = 126946.75
$ws.Range("J80").Value = 126946.75
$ws.Range("L80").Value = 126946.75
$ws.Range("N80").Value = -128942.75
$ws.Range("H83").Value = 126946.75
$ws.Range("J83").Value = 126946.75
$ws.Range("L83").Value = 380840.25
$ws.Range("N83").Value = -390824.25
$ws.Range("H97").Value = 1421.5
$ws.Range("I97").Value = 746.1667
$ws.Range("K97").Value = 746.1667
$ws.Range("M97").Value = -250.1667
$ws.Range("H102").Value = 1774.7368
$ws.Range("I102").Value = 1858.2059
$ws.Range("K102").Value = 1858.2059
$ws.Range("M102").Value = -236.2058999999999
$ws.Range("H110").Value = 2015.125
$ws.Range("I110").Value = 2179.7896
$ws.Range("J110").Value = 1389.4
$ws.Range("K110").Value = 2179.7896
$ws.Range("L110").Value = 1389.4
$ws.Range("M110").Value = -134.7896000000001
$ws.Range("N110").Value = -5479.4
$ws.Range("H116").Value = 2159.4722
$ws.Range("I116").Value = 1742.6
$ws.Range("J116").Value = 2680.5625
$ws.Range("K116").Value = 1742.6
$ws.Range("L116").Value = 2680.5625
$ws.Range("M116").Value = 551.4000000000001
$ws.Range("N116").Value = -7268.5625
$ws.Range("H132").Value = 521073.84
$ws.Range("I132").Value = 832939.9399999999
$ws.Range("J132").Value = 7412.0586
$ws.Range("K132").Value = 2498819.82
$ws.Range("L132").Value = 22236.1758
$ws.Range("M132").Value = -2496289.82
$ws.Range("N132").Value = -27296.1758

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2159.4722
$ws.Range("I3").Value = 1742.6
$ws.Range("J3").Value = 2680.5625
$ws.Range("K3").Value = 1742.6
$ws.Range("L3").Value = 2680.5625
$ws.Range("M3").Value = -1628.6
$ws.Range("N3").Value = -2908.5625
$ws.Range("H4").Value = 391.7
$ws.Range("I4").Value = 367.14285
$ws.Range("J4").Value = 449
$ws.Range("K4").Value = 367.14285
$ws.Range("L4").Value = 449
$ws.Range("M4").Value = -252.14285
$ws.Range("N4").Value = -679
$ws.Range("H9").Value = 24999.5
$ws.Range("J9").Value = 24999.5
$ws.Range("L9").Value = 24999.5
$ws.Range("N9").Value = -25335.5
$ws.Range("H20").Value = 1601.8518
$ws.Range("I20").Value = 1800.8667
$ws.Range("K20").Value = 1800.8667
$ws.Range("M20").Value = -1553.8667
$ws.Range("H86").Value = 2292.3333
$ws.Range("I86").Value = 2331.6667
$ws.Range("J86").Value = 2253
$ws.Range("K86").Value = 2331.6667
$ws.Range("L86").Value = 2253
$ws.Range("M86").Value = -1208.6667
$ws.Range("N86").Value = -4499
$ws.Range("H89").Value = 2292.3333
$ws.Range("I89").Value = 2331.6667
$ws.Range("J89").Value = 2253
$ws.Range("K89").Value = 11658.3335
$ws.Range("L89").Value = 11265
$ws.Range("M89").Value = -6042.333500000001
$ws.Range("N89").Value = -22497
$ws.Range("H94").Value = 1262.3
$ws.Range("I94").Value = 1360.8572
$ws.Range("J94").Value = 1032.3334
$ws.Range("K94").Value = 1360.8572
$ws.Range("L94").Value = 1032.3334
$ws.Range("M94").Value = -909.8571999999999
$ws.Range("N94").Value = -1934.3334
$ws.Range("H99").Value = 2023.3448
$ws.Range("I99").Value = 1572.762
$ws.Range("K99").Value = 1572.762
$ws.Range("M99").Value = -74.76199999999994
$ws.Range("H105").Value = 4110.9414
$ws.Range("I105").Value = 3859.2
$ws.Range("K105").Value = 3859.2
$ws.Range("M105").Value = -2112.2
$ws.Range("H134").Value = 1260369.8
$ws.Range("I134").Value = 1989819.4
$ws.Range("J134").Value = 9884.857
$ws.Range("K134").Value = 5969458.199999999
$ws.Range("L134").Value = 29654.571
$ws.Range("M134").Value = -5966923.199999999
$ws.Range("N134").Value = -34724.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 155.625
$ws.Range("I7").Value = 73.15385000000001
$ws.Range("J7").Value = 513
$ws.Range("K7").Value = 73.15385000000001
$ws.Range("L7").Value = 513
$ws.Range("M7").Value = 39.84614999999999
$ws.Range("N7").Value = -739
$ws.Range("H22").Value = 732.06665
$ws.Range("I22").Value = 715.1667
$ws.Range("J22").Value = 799.6667
$ws.Range("K22").Value = 715.1667
$ws.Range("L22").Value = 799.6667
$ws.Range("M22").Value = -365.1667
$ws.Range("N22").Value = -1499.6667
$ws.Range("H31").Value = 12124.8125
$ws.Range("I31").Value = 2402.3333
$ws.Range("J31").Value = 14368.462
$ws.Range("K31").Value = 2402.3333
$ws.Range("L31").Value = 14368.462
$ws.Range("M31").Value = -2107.3333
$ws.Range("N31").Value = -14958.462
$ws.Range("H34").Value = 12124.8125
$ws.Range("I34").Value = 2402.3333
$ws.Range("J34").Value = 14368.462
$ws.Range("K34").Value = 2402.3333
$ws.Range("L34").Value = 14368.462
$ws.Range("M34").Value = -2200.3333
$ws.Range("N34").Value = -14772.462
$ws.Range("H58").Value = 377745.44
$ws.Range("I58").Value = 478385.8
$ws.Range("K58").Value = 478385.8
$ws.Range("M58").Value = -478182.8
$ws.Range("H86").Value = 2379.5881
$ws.Range("I86").Value = 2226.6155
$ws.Range("K86").Value = 2226.6155
$ws.Range("M86").Value = -1103.6155
$ws.Range("H89").Value = 2379.5881
$ws.Range("I89").Value = 2226.6155
$ws.Range("K89").Value = 11133.0775
$ws.Range("M89").Value = -5517.077499999999
$ws.Range("H94").Value = 2359
$ws.Range("I94").Value = 1370.3334
$ws.Range("J94").Value = 2853.3333
$ws.Range("K94").Value = 1370.3334
$ws.Range("L94").Value = 2853.3333
$ws.Range("M94").Value = -919.3334
$ws.Range("N94").Value = -3755.3333
$ws.Range("H99").Value = 6053
$ws.Range("I99").Value = 4377.75
$ws.Range("J99").Value = 7728.25
$ws.Range("K99").Value = 4377.75
$ws.Range("L99").Value = 7728.25
$ws.Range("M99").Value = -2879.75
$ws.Range("N99").Value = -10724.25
$ws.Range("H105").Value = 21570.295
$ws.Range("I105").Value = 22862.188
$ws.Range("K105").Value = 22862.188
$ws.Range("M105").Value = -21115.188
$ws.Range("H107").Value = 1310.3914
$ws.Range("I107").Value = 913.64703
$ws.Range("J107").Value = 2434.5
$ws.Range("K107").Value = 913.64703
$ws.Range("L107").Value = 2434.5
$ws.Range("M107").Value = 1006.35297
$ws.Range("N107").Value = -6274.5
$ws.Range("H126").Value = 6053
$ws.Range("I126").Value = 4377.75
$ws.Range("J126").Value = 7728.25
$ws.Range("K126").Value = 13133.25
$ws.Range("L126").Value = 23184.75
$ws.Range("M126").Value = -10663.25
$ws.Range("N126").Value = -28124.75
$ws.Range("H132").Value = 27823616
$ws.Range("I132").Value = 119537.664
$ws.Range("K132").Value = 358612.992
$ws.Range("M132").Value = -356082.992
$ws.Range("H134").Value = 4935.1904
$ws.Range("I134").Value = 1750.8462
$ws.Range("K134").Value = 5252.5386
$ws.Range("M134").Value = -2717.5386
$ws.Range("H135").Value = 98333.164
$ws.Range("J135").Value = 98333.164
$ws.Range("L135").Value = 98333.164
$ws.Range("N135").Value = -108473.164
$ws.Range("H136").Value = 377745.44
$ws.Range("I136").Value = 478385.8
$ws.Range("K136").Value = 1435157.4
$ws.Range("M136").Value = -1432607.4
$ws.Range("H140").Value = 99999.86
$ws.Range("J140").Value = 99999.86
$ws.Range("L140").Value = 99999.86
$ws.Range("N140").Value = -110359.86

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2800.6316
$ws.Range("I5").Value = 567.5333000000001
$ws.Range("K5").Value = 1702.5999
$ws.Range("M5").Value = -1590.5999
$ws.Range("H6").Value = 81.15385000000001
$ws.Range("I6").Value = 21.444445
$ws.Range("J6").Value = 215.5
$ws.Range("K6").Value = 64.33333500000001
$ws.Range("L6").Value = 646.5
$ws.Range("M6").Value = 48.66666499999999
$ws.Range("N6").Value = -872.5
$ws.Range("H39").Value = 19004
$ws.Range("J39").Value = 19004
$ws.Range("L39").Value = 57012
$ws.Range("N39").Value = -57600
$ws.Range("H97").Value = 455.6
$ws.Range("I97").Value = 344.5
$ws.Range("K97").Value = 1033.5
$ws.Range("M97").Value = -537.5
$ws.Range("H98").Value = 999
$ws.Range("J98").Value = 999
$ws.Range("L98").Value = 2997
$ws.Range("N98").Value = -5993
$ws.Range("H109").Value = 2367
$ws.Range("I109").Value = 1380.5
$ws.Range("J109").Value = 6313
$ws.Range("K109").Value = 4141.5
$ws.Range("L109").Value = 18939
$ws.Range("M109").Value = -3101.5
$ws.Range("N109").Value = -21019
$ws.Range("H135").Value = 2800.6316
$ws.Range("I135").Value = 567.5333000000001
$ws.Range("K135").Value = 5107.7997
$ws.Range("M135").Value = -2572.7997
$ws.Range("H137").Value = 3995.4285
$ws.Range("J137").Value = 4995
$ws.Range("L137").Value = 14985
$ws.Range("N137").Value = -25185

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 204.94118
$ws.Range("I2").Value = 72.40000000000001
$ws.Range("J2").Value = 394.2857
$ws.Range("K2").Value = 72.40000000000001
$ws.Range("L2").Value = 394.2857
$ws.Range("M2").Value = 40.59999999999999
$ws.Range("N2").Value = -620.2857
$ws.Range("H80").Value = 317056
$ws.Range("I80").Value = 459036.1
$ws.Range("K80").Value = 459036.1
$ws.Range("M80").Value = -458038.1
$ws.Range("H83").Value = 317056
$ws.Range("I83").Value = 459036.1
$ws.Range("K83").Value = 2295180.5
$ws.Range("M83").Value = -2290188.5
$ws.Range("H102").Value = 2382.2856
$ws.Range("I102").Value = 2307.0454
$ws.Range("K102").Value = 2307.0454
$ws.Range("M102").Value = -685.0454
$ws.Range("H122").Value = 4019.4
$ws.Range("I122").Value = 2808.476
$ws.Range("K122").Value = 8425.428
$ws.Range("M122").Value = -5975.428
$ws.Range("H132").Value = 2407.7812
$ws.Range("I132").Value = 1822.5385
$ws.Range("J132").Value = 4943.8335
$ws.Range("K132").Value = 5467.6155
$ws.Range("L132").Value = 14831.5005
$ws.Range("M132").Value = -2937.6155
$ws.Range("N132").Value = -19891.5005
$ws.Range("H135").Value = 104998.75
$ws.Range("J135").Value = 104998.75
$ws.Range("L135").Value = 104998.75
$ws.Range("N135").Value = -115138.75
$ws.Range("H140").Value = 106747.125
$ws.Range("J140").Value = 106747.125
$ws.Range("L140").Value = 106747.125
$ws.Range("N140").Value = -117107.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4822
$ws.Range("J16").Value = 7500
$ws.Range("L16").Value = 7500
$ws.Range("N16").Value = -7840
$ws.Range("H22").Value = 30051.445
$ws.Range("I22").Value = 84497.914
$ws.Range("J22").Value = 2828.2083
$ws.Range("K22").Value = 84497.914
$ws.Range("L22").Value = 2828.2083
$ws.Range("M22").Value = -84202.914
$ws.Range("N22").Value = -3418.2083
$ws.Range("H27").Value = 30051.445
$ws.Range("I27").Value = 84497.914
$ws.Range("J27").Value = 2828.2083
$ws.Range("K27").Value = 84497.914
$ws.Range("L27").Value = 2828.2083
$ws.Range("M27").Value = -84390.914
$ws.Range("N27").Value = -3042.2083
$ws.Range("H40").Value = 1543.4615
$ws.Range("I40").Value = 1543.4615
$ws.Range("K40").Value = 1543.4615
$ws.Range("M40").Value = -1407.4615
$ws.Range("H132").Value = 668827.8
$ws.Range("I132").Value = 991158.5
$ws.Range("J132").Value = 5205.8237
$ws.Range("K132").Value = 2973475.5
$ws.Range("L132").Value = 15617.4711
$ws.Range("M132").Value = -2970945.5
$ws.Range("N132").Value = -20677.4711
$ws.Range("H136").Value = 8687.625
$ws.Range("I136").Value = 3832.3333
$ws.Range("K136").Value = 11496.9999
$ws.Range("M136").Value = -8946.999899999999
$ws.Range("H139").Value = 98810
$ws.Range("J139").Value = 98810
$ws.Range("L139").Value = 98810
$ws.Range("N139").Value = -109090
$ws.Range("H140").Value = 114870.4
$ws.Range("J140").Value = 114870.4
$ws.Range("L140").Value = 114870.4
$ws.Range("N140").Value = -125230.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 44444
$ws.Range("J16").Value = 44444
$ws.Range("L16").Value = 44444
$ws.Range("N16").Value = -45028
$ws.Range("H46").Value = 77775.22
$ws.Range("J46").Value = 78497.125
$ws.Range("L46").Value = 78497.125
$ws.Range("N46").Value = -78959.125
$ws.Range("H54").Value = 123950
$ws.Range("J54").Value = 123950
$ws.Range("L54").Value = 123950
$ws.Range("N54").Value = -124990
$ws.Range("H81").Value = 2455.4443
$ws.Range("I81").Value = 2455.4443
$ws.Range("K81").Value = 4910.8886
$ws.Range("M81").Value = -3849.8886
$ws.Range("H84").Value = 2455.4443
$ws.Range("I84").Value = 2455.4443
$ws.Range("K84").Value = 24554.443
$ws.Range("M84").Value = -19250.443
$ws.Range("H100").Value = 1719.1904
$ws.Range("I100").Value = 2118.4614
$ws.Range("K100").Value = 4236.9228
$ws.Range("M100").Value = -3695.9228
$ws.Range("H107").Value = 1708.7368
$ws.Range("I107").Value = 1321.5294
$ws.Range("K107").Value = 3964.5882
$ws.Range("M107").Value = -2044.5882
$ws.Range("H122").Value = 4037.5293
$ws.Range("I122").Value = 3664.5386
$ws.Range("K122").Value = 10993.6158
$ws.Range("M122").Value = -8543.6158
$ws.Range("H126").Value = 3663.3333
$ws.Range("I126").Value = 3196.6
$ws.Range("J126").Value = 5997
$ws.Range("K126").Value = 9589.799999999999
$ws.Range("L126").Value = 17991
$ws.Range("M126").Value = -7119.799999999999
$ws.Range("N126").Value = -22931
$ws.Range("H134").Value = 77775.22
$ws.Range("J134").Value = 78497.125
$ws.Range("L134").Value = 235491.375
$ws.Range("N134").Value = -240561.375
$ws.Range("H136").Value = 14072882
$ws.Range("I136").Value = 22349014
$ws.Range("K136").Value = 67047042
$ws.Range("M136").Value = -67044492
